$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cohort creation now estimates schedule time loss, so the previously
# zero/under-estimated "# of Students" figures for a few cohort rows are
# updated to reflect the corrected calculation.
$ws.Range("C2").Value = 40
$ws.Range("C19").Value = 20
$ws.Range("C23").Value = 10
$ws.Range("C24").Value = 10

# Reflect the saved session's cursor position/selection.
$ws.Range("O21").Select()
